$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "electrode_name"
$ws.Range("B1").Value = "recording_scale"

$ws.Range("M7").Select()
